$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "FT232500M2TTWGFH"
